# Auto-generated edit script for violent-crime-full-year workbook
# Commit: Add data for 2023-12-18
# Updates column J (year 2023 totals) on the affected worksheets to reflect
# crime counts after adding one additional day (2023-12-18) of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 7392
$ws.Cells.Item(3, 10).Value = 7789
$ws.Cells.Item(4, 10).Value = 1695
$ws.Cells.Item(5, 10).Value = 608
$ws.Cells.Item(6, 10).Value = 10626
$ws.Cells.Item(7, 10).Value = 28110

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(3, 10).Value = 10
$ws.Cells.Item(5, 10).Value = 84
$ws.Cells.Item(6, 10).Value = 215
$ws.Cells.Item(7, 10).Value = 804
$ws.Cells.Item(8, 10).Value = 1770
$ws.Cells.Item(10, 10).Value = 203
$ws.Cells.Item(11, 10).Value = 500
$ws.Cells.Item(15, 10).Value = 347
$ws.Cells.Item(16, 10).Value = 112
$ws.Cells.Item(18, 10).Value = 227
$ws.Cells.Item(19, 10).Value = 811
$ws.Cells.Item(20, 10).Value = 604
$ws.Cells.Item(24, 10).Value = 95
$ws.Cells.Item(27, 10).Value = 168
$ws.Cells.Item(29, 10).Value = 1496
$ws.Cells.Item(33, 10).Value = 1275
$ws.Cells.Item(36, 10).Value = 381
$ws.Cells.Item(37, 10).Value = 860
$ws.Cells.Item(40, 10).Value = 63
$ws.Cells.Item(42, 10).Value = 1195
$ws.Cells.Item(44, 10).Value = 221
$ws.Cells.Item(45, 10).Value = 40
$ws.Cells.Item(47, 10).Value = 204
$ws.Cells.Item(48, 10).Value = 314
$ws.Cells.Item(51, 10).Value = 352
$ws.Cells.Item(52, 10).Value = 714
$ws.Cells.Item(53, 10).Value = 423
$ws.Cells.Item(54, 10).Value = 557
$ws.Cells.Item(55, 10).Value = 442
$ws.Cells.Item(63, 10).Value = 95
$ws.Cells.Item(64, 10).Value = 188
$ws.Cells.Item(65, 10).Value = 707
$ws.Cells.Item(67, 10).Value = 1030
$ws.Cells.Item(73, 10).Value = 274
$ws.Cells.Item(75, 10).Value = 85
$ws.Cells.Item(76, 10).Value = 399
$ws.Cells.Item(78, 10).Value = 326
$ws.Cells.Item(79, 10).Value = 768
$ws.Cells.Item(83, 10).Value = 565
$ws.Cells.Item(85, 10).Value = 1155
$ws.Cells.Item(87, 10).Value = 95
$ws.Cells.Item(88, 10).Value = 298
$ws.Cells.Item(89, 10).Value = 349
$ws.Cells.Item(91, 10).Value = 322
$ws.Cells.Item(92, 10).Value = 91
$ws.Cells.Item(95, 10).Value = 404
$ws.Cells.Item(96, 10).Value = 313
$ws.Cells.Item(97, 10).Value = 255
$ws.Cells.Item(98, 10).Value = 208
$ws.Cells.Item(99, 10).Value = 427
$ws.Cells.Item(101, 10).Value = 28110

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 10).Value = 81
$ws.Cells.Item(6, 10).Value = 118
$ws.Cells.Item(7, 10).Value = 313

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 254
$ws.Cells.Item(3, 10).Value = 240
$ws.Cells.Item(7, 10).Value = 804

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 10).Value = 86
$ws.Cells.Item(6, 10).Value = 237
$ws.Cells.Item(7, 10).Value = 500

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 10).Value = 106
$ws.Cells.Item(7, 10).Value = 349

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 307
$ws.Cells.Item(3, 10).Value = 418
$ws.Cells.Item(6, 10).Value = 329
$ws.Cells.Item(7, 10).Value = 1155

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(7, 10).Value = 714

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 10).Value = 276
$ws.Cells.Item(7, 10).Value = 423

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 10).Value = 512
$ws.Cells.Item(6, 10).Value = 656
$ws.Cells.Item(7, 10).Value = 1770

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 210
$ws.Cells.Item(6, 10).Value = 151
$ws.Cells.Item(7, 10).Value = 565

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 289
$ws.Cells.Item(3, 10).Value = 424
$ws.Cells.Item(5, 10).Value = 54
$ws.Cells.Item(6, 10).Value = 453
$ws.Cells.Item(7, 10).Value = 1275

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 10).Value = 142
$ws.Cells.Item(3, 10).Value = 147
$ws.Cells.Item(7, 10).Value = 404

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 260
$ws.Cells.Item(7, 10).Value = 860

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 10).Value = 205
$ws.Cells.Item(7, 10).Value = 707

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 117
$ws.Cells.Item(7, 10).Value = 427

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 389
$ws.Cells.Item(6, 10).Value = 284
$ws.Cells.Item(7, 10).Value = 1030

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 10).Value = 138
$ws.Cells.Item(3, 10).Value = 113
$ws.Cells.Item(6, 10).Value = 257
$ws.Cells.Item(7, 10).Value = 557

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 455
$ws.Cells.Item(3, 10).Value = 527
$ws.Cells.Item(6, 10).Value = 379
$ws.Cells.Item(7, 10).Value = 1496

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 10).Value = 53
$ws.Cells.Item(3, 10).Value = 61
$ws.Cells.Item(6, 10).Value = 150
$ws.Cells.Item(7, 10).Value = 314

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 10).Value = 199
$ws.Cells.Item(3, 10).Value = 230
$ws.Cells.Item(4, 10).Value = 38
$ws.Cells.Item(6, 10).Value = 314
$ws.Cells.Item(7, 10).Value = 811

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 10).Value = 89
$ws.Cells.Item(7, 10).Value = 221

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 10).Value = 206
$ws.Cells.Item(7, 10).Value = 399

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 10).Value = 77
$ws.Cells.Item(7, 10).Value = 215

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 10).Value = 633
$ws.Cells.Item(7, 10).Value = 1195

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 10).Value = 116
$ws.Cells.Item(7, 10).Value = 203

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 326

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 10).Value = 85
$ws.Cells.Item(7, 10).Value = 442

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(6, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 95

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 10).Value = 84
$ws.Cells.Item(6, 10).Value = 84
$ws.Cells.Item(7, 10).Value = 322

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 10).Value = 257
$ws.Cells.Item(7, 10).Value = 768

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(6, 10).Value = 67
$ws.Cells.Item(7, 10).Value = 188

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 10).Value = 198
$ws.Cells.Item(4, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 604

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 10).Value = 106
$ws.Cells.Item(7, 10).Value = 227

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(6, 10).Value = 113
$ws.Cells.Item(7, 10).Value = 381

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 204

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 10).Value = 71
$ws.Cells.Item(7, 10).Value = 347

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(2, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 208

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 10).Value = 86
$ws.Cells.Item(3, 10).Value = 70
$ws.Cells.Item(7, 10).Value = 274

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 255

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(6, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 91

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(6, 10).Value = 158
$ws.Cells.Item(7, 10).Value = 298

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(3, 10).Value = 16
$ws.Cells.Item(7, 10).Value = 84

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 10).Value = 44
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 168

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(3, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 85

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 10).Value = 74
$ws.Cells.Item(3, 10).Value = 92
$ws.Cells.Item(7, 10).Value = 352

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Cells.Item(3, 10).Value = 10
$ws.Cells.Item(7, 10).Value = 40

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(6, 10).Value = 10

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(3, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 63

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 10).Value = 62
$ws.Cells.Item(7, 10).Value = 95

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 112
